# Recalibrate retirements per unit net loss parameter
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRpUNL")

# Remove the gray "calibration highlight" shading that was covering column B
# (rows 2-25); the workbook no longer needs it so the style/fill falls out
# of use entirely.
$ws.Range("B2:B25").ClearFormats()

# Recalibrate the base retirement-per-unit-net-loss parameter.
$ws.Range("B2").Value = 0.015

# All of the other plant types that aren't pinned at 0 (the "very low"
# retirement-fraction rows) should simply reference the recalibrated base
# value instead of repeating the old hard-coded 0.03 constant.
$ws.Range("B3").Formula = '=$B$2'
$ws.Range("B4").Formula = '=$B$2'
$ws.Range("B5").Formula = '=$B$2'
$ws.Range("B7").Formula = '=$B$2'
$ws.Range("B8").Formula = '=$B$2'
$ws.Range("B13").Formula = '=$B$2'
$ws.Range("B14").Formula = '=$B$2'
$ws.Range("B15").Formula = '=$B$2'
$ws.Range("B19").Formula = '=$B$2'
$ws.Range("B20").Formula = '=$B$2'
$ws.Range("B21").Formula = '=$B$2'
$ws.Range("B22").Formula = '=$B$2'
$ws.Range("B23").Formula = '=$B$2'
$ws.Range("B24").Formula = '=$B$2'
$ws.Range("B25").Formula = '=$B$2'

# Leave the cursor on the cell the author was last looking at, and make the
# CRpUNL tab the active/selected sheet (it was "About" before).
$ws.Range("C14").Select() | Out-Null
$ws.Activate() | Out-Null
